# feat: add 2022-Q1 data
#
# The workbook's running "总计" (totals) sheet is repurposed in place to
# become the new "2022-Q1" per-fund holdings sheet (it keeps its original
# sheetId/tab position), and a fresh "总计" sheet is appended right after it
# carrying the refreshed roll-up table (new 2022-Q1 row + the previous
# totals, still in reverse-chronological order).

$wb = $excel.ActiveWorkbook

# A never-touched cell we use purely as a "blank formatting" donor, so we
# can paste-reset a cell's style back to the workbook default after typing
# a quote-prefixed value into it (Value = "'123" stores the text correctly
# but also stamps a quote-prefix style on the cell; pasting blank formats
# over it clears that stamp while leaving the stored text untouched).
$blank = $wb.Worksheets.Item(1).Cells.Item(100, 100)

# ============================================================
# Step 1: turn the existing "总计" sheet into "2022-Q1"
# ============================================================
$q1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1.Cells.ClearContents()
$q1.Name = "2022-Q1"

# Drop the two surplus rows left over from the old 6-row totals table -
# the new fund table only needs a header + 3 data rows.
$q1.Rows(5).Delete()
$q1.Rows(5).Delete()

$q1Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $q1Headers.Length; $i++) {
    $q1.Cells.Item(1, $i + 2).Value = $q1Headers[$i]
}
# B1:D1 already carry the header style inherited from the old sheet;
# extend the same formatting across the newly-used E1:H1 header cells.
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

$q1Rows = @(
    @("000586", "景顺长城中小板创业板精选股票", "2.42", "94.15", "8.39", "0.2030", 1),
    @("010706", "景顺长城景骊成长混合型证券投资基金", "1.13", "93.50", "6.37", "0.0720", 2),
    @("260115", "景顺长城中小盘混合", "0.96", "94.00", "6.55", "0.0629", 2)
)

for ($r = 0; $r -lt $q1Rows.Length; $r++) {
    $row = $q1Rows[$r]
    $excelRow = $r + 2
    # index column keeps the style inherited from the old sheet already
    $q1.Cells.Item($excelRow, 1).Value = $r
    $q1.Cells.Item($excelRow, 3).Value = $row[1]
    $q1.Cells.Item($excelRow, 8).Value = $row[6]

    # text-like numeric-looking columns (fund code / scale / position /
    # ratio / market value) must stay plain text, not be coerced to
    # numbers - enter with a quote prefix, then strip the resulting
    # quote-prefix style stamp back to the sheet's plain default.
    $q1.Cells.Item($excelRow, 2).Value = "'" + $row[0]
    $q1.Cells.Item($excelRow, 4).Value = "'" + $row[2]
    $q1.Cells.Item($excelRow, 5).Value = "'" + $row[3]
    $q1.Cells.Item($excelRow, 6).Value = "'" + $row[4]
    $q1.Cells.Item($excelRow, 7).Value = "'" + $row[5]
}
$blank.Copy()
$q1.Range("B2:B4").PasteSpecial(-4122)
$q1.Range("D2:G4").PasteSpecial(-4122)

# ============================================================
# Step 2: append a fresh "总计" sheet after "2022-Q1"
# ============================================================
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# Pull in the same header/index-column styling used throughout the
# workbook's other sheets.
$q1.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$q1.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($i = 0; $i -lt $totalHeaders.Length; $i++) {
    $total.Cells.Item(1, $i + 2).Value = $totalHeaders[$i]
}

$totalRows = @(
    @("2022-Q1", 3, 0.34),
    @("2021-Q4", 2, 0.27),
    @("2021-Q3", 3, 0.31),
    @("2021-Q2", 3, 0.34),
    @("2021-Q1", 3, 0.35),
    @("2020-Q4", 3, 0.46)
)

for ($r = 0; $r -lt $totalRows.Length; $r++) {
    $row = $totalRows[$r]
    $excelRow = $r + 2
    $total.Cells.Item($excelRow, 1).Value = $r
    $total.Cells.Item($excelRow, 2).Value = $row[0]
    $total.Cells.Item($excelRow, 3).Value = $row[1]
    $total.Cells.Item($excelRow, 4).Value = $row[2]
}
